$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shelter assignments (D/E/F columns) for a few communities.
$ws.Range("D4").Value = "Old Municipal Bldg."
$ws.Range("E4").Value = 120.948177254006
$ws.Range("F4").Value = 14.7573006861396

$ws.Range("D6").Value = "Ibayo Elementary School"
$ws.Range("E6").Value = 120.959816737558
$ws.Range("F6").Value = 14.7535649557989

$ws.Range("D9").Value = "FSS Patulo Elementary School"
$ws.Range("E9").Value = 121.027062736924
$ws.Range("F9").Value = 14.7839553140957

$ws.Range("D10").Value = "Barangay Hall Nagbalon"
$ws.Range("E10").Value = 120.950788291388
$ws.Range("F10").Value = 14.7523618894178

# Delete the "Saog" community row (row 16); "Tabing Ilog" (row 17) shifts up to row 16.
$ws.Rows.Item(16).Delete()
